# Update inputs sheet with abattoir worker mask info (by gender) and
# clarify the clearing_day description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: clarify clearing_day description (max. 36) ---
$ws.Range("D3").Value = "Day index when flock is cleared (max. 36)"

# --- Row 16: p_mask -> p_mask_male ---
$ws.Range("A16").Value = "p_mask_male"
$ws.Range("B16").Value = "Selected"
$ws.Range("C16").Value = "DOUBLE"
$ws.Range("D16").Value = "Probability of a male worker wearing mask"
$ws.Range("E16").Value = "User defined"
$ws.Range("F16").Value = "proportion"
$ws.Range("G16").Value = 1

# --- Row 17: p_wash -> p_mask_female (new meaning) ---
$ws.Range("A17").Value = "p_mask_female"
$ws.Range("B17").Value = "Selected"
$ws.Range("C17").Value = "DOUBLE"
$ws.Range("D17").Value = "Probability of a female worker wearing mask"
$ws.Range("E17").Value = "User defined"
$ws.Range("F17").Value = "proportion"
$ws.Range("G17").Value = 0

# --- Row 18: p_glove -> p_wash ---
$ws.Range("A18").Value = "p_wash"
$ws.Range("B18").Value = "Selected"
$ws.Range("C18").Value = "DOUBLE"
$ws.Range("D18").Value = "Probability of a worker washed hands afterwards"
$ws.Range("E18").Value = "User defined"
$ws.Range("F18").Value = "proportion"
$ws.Range("G18").Value = 1

# --- Row 19 (new): p_glove ---
$ws.Range("A19").Value = "p_glove"
$ws.Range("B19").Value = "Selected"
$ws.Range("C19").Value = "DOUBLE"
$ws.Range("D19").Value = "Probability of a worker wearing gloves"
$ws.Range("E19").Value = "User defined"
$ws.Range("F19").Value = "proportion"
$ws.Range("G19").Value = 1

# --- Row 20 (new): p_male ---
$ws.Range("A20").Value = "p_male"
$ws.Range("B20").Value = "Selected"
$ws.Range("C20").Value = "DOUBLE"
$ws.Range("D20").Value = "Proportion of male workers"
$ws.Range("E20").Value = "User defined"
$ws.Range("F20").Value = "proportion"
$ws.Range("G20").Value = 0.6

# --- Update selection to match final state ---
$ws.Range("K7").Select()
